# Pin Assignments.xlsx edit
# - Adds a "Channel #" / "wire color" pair of columns to both the Read
#   Operations and Write Operations tables.
# - Re-orders the G0/G1 rows in the Read Operations table.
# - Moves the Write Operations table two columns to the right (F:H -> H:L)
#   to make room for the new columns.
# - Adds a new row (WS Home Switch / Limit Switch Return) to both tables.
# - Adjusts column widths, the workbook window and the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the whole area that is being rebuilt -
# both tables are fully re-laid-out so it's simplest to wipe the old
# content/format and rewrite every cell explicitly.
$ws.Range("A1:L8").Clear()

# ---------------------------------------------------------------------
# Values (text cells)
# ---------------------------------------------------------------------
$textCells = @(
    @("A1", "Read Operations"),
    @("H1", "Write Opeatrions"),

    @("A2", "Function"),
    @("B2", "Pin #"),
    @("C2", "Pin Type"),
    @("D2", "Channel #"),
    @("E2", "wire color"),
    @("H2", "Function"),
    @("I2", "Pin #"),
    @("J2", "Pin Type"),
    @("K2", "Channel #"),
    @("L2", "wire color"),

    @("A3", "G1"),
    @("B3", "Pin 9"),
    @("C3", "DIO"),
    @("E3", "Solid Orange"),
    @("H3", "G1"),
    @("I3", "Pin 13"),
    @("J3", "DIO"),
    @("L3", "Green"),

    @("A4", "G0"),
    @("B4", "Pin 11"),
    @("C4", "DIO"),
    @("E4", "Solid Yellow"),
    @("H4", "G0"),
    @("I4", "Pin 15"),
    @("J4", "DIO"),
    @("L4", "Orange"),

    @("A5", "Temp"),
    @("B5", "Pin 1"),
    @("C5", "Analog"),
    @("H5", "Pulse"),
    @("I5", "Pin5"),
    @("J5", "DIO"),
    @("L5", "yellow/black"),

    @("A6", "Current"),
    @("B6", "Pin 3"),
    @("C6", "Analog"),
    @("E6", "Solid Green"),
    @("H6", "Direction"),
    @("I6", "Pin 3"),
    @("J6", "DIO"),
    @("L6", "orange/black"),

    @("A7", "Measure WS"),
    @("B7", "Pin 1-9"),
    @("C7", "Counter"),
    @("E7", "See CA-FC10 Datasheet"),
    @("H7", "Stop/Start CNC"),
    @("I7", "Pin 1"),
    @("J7", "DIO"),
    @("L7", "Solid Teal"),

    @("A8", "WS Home Switch"),
    @("B8", "Pin 7"),
    @("C8", "DIO"),
    @("E8", "purple"),
    @("H8", "Limit Switch Return"),
    @("I8", "Pin 4"),
    @("J8", "GND"),
    @("K8", "n/a"),
    @("L8", "black/white")
)

foreach ($pair in $textCells) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# ---------------------------------------------------------------------
# Values (numeric "Channel #" cells)
# ---------------------------------------------------------------------
$numberCells = @(
    @("D3", 19),
    @("D4", 18),
    @("D8", 20),
    @("K3", 17),
    @("K4", 16),
    @("K5", 21),
    @("K6", 22),
    @("K7", 23)
)

foreach ($pair in $numberCells) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# D5, D6, D7 stay blank but keep the bordered/formatted look of the rest
# of the "Channel #" column - handled by the border formatting below.

# ---------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------
# NOTE: this engine's Range object only honours the FIRST area of a
# multi-area (comma-separated) reference, so each table (A:E vs H:L) is
# formatted with its own statement rather than a single unioned Range.

# Title cells: bold, larger font, no border (matches A1 style).
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 14
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Font.Size = 14

# Column header rows: bold, bordered.
$ws.Range("A2:E2").Font.Bold = $true
$ws.Range("A2:E2").Font.Size = 12
$ws.Range("A2:E2").Borders.LineStyle = 1
$ws.Range("H2:L2").Font.Bold = $true
$ws.Range("H2:L2").Font.Size = 12
$ws.Range("H2:L2").Borders.LineStyle = 1

# Data cells: regular font, bordered. Includes the blank D5:D7 cells so
# they keep the table's grid lines even though they hold no value.
$ws.Range("A3:E8").Font.Bold = $false
$ws.Range("A3:E8").Font.Size = 12
$ws.Range("A3:E8").Borders.LineStyle = 1
$ws.Range("H3:L8").Font.Bold = $false
$ws.Range("H3:L8").Font.Size = 12
$ws.Range("H3:L8").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 14.67   # A  -> raw width 15.5
$ws.Columns(5).ColumnWidth = 11.17   # E  -> raw width 12
$ws.Columns(6).ColumnWidth = 12.67   # F  -> raw width 13.5 (unchanged)
$ws.Columns(8).ColumnWidth = 17.5    # H  -> raw width ~18.33

# ---------------------------------------------------------------------
# Window / selection
# ---------------------------------------------------------------------
# Best-effort: move/resize the saved window to the new position recorded
# in the workbook (xWindow=0, windowWidth=25600). Not all hosts persist
# window placement back to the file, but setting it is harmless.
$aw = $excel.ActiveWindow
$aw.Left = 0
$aw.Top = 0
$aw.Width = 25600
$aw.Height = 15480

$null = $ws.Range("E5").Select()
